# Adds "Page classes for requests" data: a new API_DATA_FILE column (G) and
# fills in the data file / data-file-column values for the API test rows
# (16-19), matching the commit "Added  Page classes for requests".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the F (TC_DATAFILE) / G (new API_DATA_FILE) columns for the
#     four API test-case rows. F16/G16 already exist (empty) in the sheet,
#     while F17:G19 are brand-new cells. ---
$ws.Range("F16").Value = "data.xlsx"
$ws.Range("G16").Value = "API_DATA_FILE"

$ws.Range("F17").Value = "data.xlsx"
$ws.Range("G17").Value = "API_DATA_FILE"

$ws.Range("F18").Value = "data.xlsx"
$ws.Range("G18").Value = "API_DATA_FILE"

$ws.Range("F19").Value = "data.xlsx"
$ws.Range("G19").Value = "API_DATA_FILE"

# --- F16:G16 already carry the correct (non-default) cell style. Copy that
#     formatting down onto the newly created F17:G19 cells so they match. ---
$ws.Range("F16:G16").Copy() | Out-Null
$ws.Range("F17:G19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column G needs to be wide enough to show "API_DATA_FILE" in full. ---
$ws.Columns("G").ColumnWidth = 16.83

# --- Reflect the resulting selection / scroll position used while editing. ---
$win = $excel.ActiveWindow
$ws.Range("G16:G19").Select() | Out-Null
$win.ScrollColumn = 2
$win.ScrollRow = 5
